$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Return_with_prediction" (G) values for rows 2-85, sourced from a re-run
# of the upstream prediction pipeline (final grouped files updated).
$gArr = New-Object 'object[,]' 84,1
$gArr[0,0] = 0.06322124964038238
$gArr[1,0] = 0.06149556876201699
$gArr[2,0] = -0.02316033021422086
$gArr[3,0] = -0.01717093012244481
$gArr[4,0] = -0.01019305773856364
$gArr[5,0] = -0.005986732640786203
$gArr[6,0] = 0.001822285315546187
$gArr[7,0] = 0.007659614218980378
$gArr[8,0] = -0.06640624874590992
$gArr[9,0] = -0.06384208373839481
$gArr[10,0] = -0.4098374533384486
$gArr[11,0] = -0.4070294283374353
$gArr[12,0] = -0.01433968120879874
$gArr[13,0] = -0.007934591074170187
$gArr[14,0] = 0.1310771767301388
$gArr[15,0] = 0.1474079812310827
$gArr[16,0] = 0.1248803713104605
$gArr[17,0] = 0.1275621745818734
$gArr[18,0] = 0.08680885261731229
$gArr[19,0] = 0.08872760067437924
$gArr[20,0] = -0.09386549115228399
$gArr[21,0] = -0.09673340381209043
$gArr[22,0] = 0.155164777170661
$gArr[23,0] = 0.1667537146473687
$gArr[24,0] = 0.08465395912220915
$gArr[25,0] = 0.08601341179831895
$gArr[26,0] = -0.1450922643496933
$gArr[27,0] = -0.1413533335968777
$gArr[28,0] = 0.03864033859799502
$gArr[29,0] = 0.04328705156343251
$gArr[30,0] = 0.1106710050374117
$gArr[31,0] = 0.1257797355280407
$gArr[32,0] = -0.01737526707036435
$gArr[33,0] = -0.01568144296554796
$gArr[34,0] = 0.03718370058490481
$gArr[35,0] = 0.03303677690914621
$gArr[36,0] = 0.1089488383339371
$gArr[37,0] = 0.1044394604154137
$gArr[38,0] = 0.03402023921723832
$gArr[39,0] = 0.03377146847379955
$gArr[40,0] = 0.1181377575560425
$gArr[41,0] = 0.1156779080988015
$gArr[42,0] = 0.03692370174260107
$gArr[43,0] = 0.03573701523118656
$gArr[44,0] = 0.05882937079726518
$gArr[45,0] = 0.06298862812376178
$gArr[46,0] = 0.04534870153724891
$gArr[47,0] = 0.04345716113199896
$gArr[48,0] = 0.023924805248048
$gArr[49,0] = 0.03115338102962048
$gArr[50,0] = -0.08069709786939101
$gArr[51,0] = -0.08142060475779567
$gArr[52,0] = 0.04488295905481981
$gArr[53,0] = 0.05540563095334069
$gArr[54,0] = 0.04634566333794014
$gArr[55,0] = 0.03833916300773737
$gArr[56,0] = 0.05538559615160537
$gArr[57,0] = 0.05711399024295274
$gArr[58,0] = 0.03268136788416055
$gArr[59,0] = 0.02824109910075292
$gArr[60,0] = 0.06229342704472279
$gArr[61,0] = 0.05740697711788335
$gArr[62,0] = 0.02692331558567181
$gArr[63,0] = 0.03103999355822355
$gArr[64,0] = 0.08034480704990898
$gArr[65,0] = 0.07847816782147224
$gArr[66,0] = -0.01870629457997228
$gArr[67,0] = -0.0250679136586608
$gArr[68,0] = 0.07197977464371674
$gArr[69,0] = 0.07862772606337291
$gArr[70,0] = -0.1504128083518064
$gArr[71,0] = -0.1567489308245939
$gArr[72,0] = 0.1416004383461187
$gArr[73,0] = 0.15487489122891
$gArr[74,0] = -0.008653971703748556
$gArr[75,0] = -0.001198252819052694
$gArr[76,0] = 0.09217150773229134
$gArr[77,0] = 0.09678283561372647
$gArr[78,0] = -0.223385308031304
$gArr[79,0] = -0.212373596833578
$gArr[80,0] = 0.1659680119698488
$gArr[81,0] = 0.1827395951155636
$gArr[82,0] = 0.1119420978731328
$gArr[83,0] = 0.1131033147676417
$ws.Range("G2:G85").Value = $gArr

# Updated "return_pct_change" (H) values recomputed from the new G values.
$hArr = New-Object 'object[,]' 84,1
$hArr[0,0] = -1.681047665113891
$hArr[1,0] = 9.40901520530412
$hArr[2,0] = -14.02509423193054
$hArr[3,0] = -50.20936911107328
$hArr[4,0] = 8.938176405628804
$hArr[5,0] = -5.702889911296169
$hArr[6,0] = 132.0308157183518
$hArr[7,0] = 239.4394739887549
$hArr[8,0] = -5.370358397192072
$hArr[9,0] = 0.3682450648548042
$hArr[10,0] = -3.860598377275327
$hArr[11,0] = -3.802473378439086
$hArr[12,0] = -74.99658366431328
$hArr[13,0] = 82.49156447332913
$hArr[14,0] = -4.1480710436375
$hArr[15,0] = 5.691049241276977
$hArr[16,0] = 6.129237133013037
$hArr[17,0] = -0.8600967859713871
$hArr[18,0] = -2.170268833280169
$hArr[19,0] = 1.879155058683411
$hArr[20,0] = -0.4089738805019414
$hArr[21,0] = 4.644500439095458
$hArr[22,0] = -3.67884118565554
$hArr[23,0] = -2.252118075525861
$hArr[24,0] = -6.618406345897593
$hArr[25,0] = 0.06010291898858405
$hArr[26,0] = -5.413981665260521
$hArr[27,0] = -1.111280925923375
$hArr[28,0] = -25.71347762193184
$hArr[29,0] = -1.21487608891768
$hArr[30,0] = 1.809569665042512
$hArr[31,0] = 1.357460395584268
$hArr[32,0] = -11.2629915034057
$hArr[33,0] = 6.308342282892266
$hArr[34,0] = 1.133603744414004
$hArr[35,0] = -7.422792265650033
$hArr[36,0] = 8.620555109197568
$hArr[37,0] = 7.227355610754564
$hArr[38,0] = 0.9873641196946753
$hArr[39,0] = 4.819382014001032
$hArr[40,0] = -2.289122100027037
$hArr[41,0] = -9.474896692646407
$hArr[42,0] = -6.914973206275338
$hArr[43,0] = 14.68009270697771
$hArr[44,0] = 3.906388202441836
$hArr[45,0] = 7.36053813149453
$hArr[46,0] = -7.935106841237132
$hArr[47,0] = -4.631754420016433
$hArr[48,0] = -9.671644335025631
$hArr[49,0] = 11.19905518581931
$hArr[50,0] = 7.163233341524585
$hArr[51,0] = -1.501445842311301
$hArr[52,0] = -10.28199476346466
$hArr[53,0] = -1.562449919231653
$hArr[54,0] = -6.243401452102393
$hArr[55,0] = 0.9350920298841494
$hArr[56,0] = -3.867904369150103
$hArr[57,0] = 0.1765920311696
$hArr[58,0] = 19.01445177449104
$hArr[59,0] = 5.780522441868657
$hArr[60,0] = -0.2462947100033555
$hArr[61,0] = -10.13685292391624
$hArr[62,0] = -2.94286921255558
$hArr[63,0] = -12.38385174550153
$hArr[64,0] = 3.426257581250673
$hArr[65,0] = -0.4978587317903787
$hArr[66,0] = 13.97701232480286
$hArr[67,0] = -30.9623572379347
$hArr[68,0] = -0.008457474214225375
$hArr[69,0] = -0.9988685676320732
$hArr[70,0] = 2.10143377413362
$hArr[71,0] = -2.399714306467275
$hArr[72,0] = -5.862786056118964
$hArr[73,0] = 2.942934813564358
$hArr[74,0] = -734.9533842575183
$hArr[75,0] = 45.7285453007581
$hArr[76,0] = 2.460316097167885
$hArr[77,0] = -0.1223193417386144
$hArr[78,0] = -3.212046383790419
$hArr[79,0] = 0.3473675636450159
$hArr[80,0] = -0.9841780120601649
$hArr[81,0] = 3.812249979353979
$hArr[82,0] = 5.476695987479776
$hArr[83,0] = 8.171963812296834
$ws.Range("H2:H85").Value = $hArr

# Updated "mean_return_pct_change" (I2) = mean of the new H column.
$ws.Range("I2").Value = -5.000064430856385
